# Update cryptos list values (price/volume columns) on the active sheet to the
# latest scraped snapshot. Leading apostrophes force text interpretation so
# numeric-looking price strings (e.g. "7.20", "226.58") keep their exact
# digits/trailing zeros instead of being auto-coerced into Numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.166.69"
$ws.Range("E2").Value = "'  +1.29%  "
$ws.Range("D3").Value = "'1.787.90"
$ws.Range("E3").Value = "'  +1.31%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'226.58"
$ws.Range("D6").Value = "'0.547"
$ws.Range("E6").Value = "'  +0.69%  "
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("E8").Value = "'  -0.51%  "
$ws.Range("E9").Value = "'  +1.69%  "
$ws.Range("E10").Value = "'  +0.42%  "
$ws.Range("E11").Value = "'  +1.01%  "
$ws.Range("D12").Value = "'2.046.29"
$ws.Range("E12").Value = "'  +1.32%  "
$ws.Range("E13").Value = "'  -1.58%  "
$ws.Range("D14").Value = "'1.790.16"
$ws.Range("E14").Value = "'  +1.25%  "
$ws.Range("E15").Value = "'  +2.30%  "
$ws.Range("D16").Value = "'34.112.41"
$ws.Range("E16").Value = "'  +1.16%  "
$ws.Range("E17").Value = "'  +1.17%  "
$ws.Range("E18").Value = "'  +2.57%  "
$ws.Range("D19").Value = "'247.14"
$ws.Range("E19").Value = "'  +3.99%  "
$ws.Range("E20").Value = "'  +0.76%  "
$ws.Range("D21").Value = "'10.92"
$ws.Range("E21").Value = "'  +3.29%  "
$ws.Range("E22").Value = "'  -0.10%  "
$ws.Range("E23").Value = "'  +1.50%  "
$ws.Range("E24").Value = "'  +0.53%  "
$ws.Range("D25").Value = "'161.14"
$ws.Range("E25").Value = "'  +1.31%  "
$ws.Range("D26").Value = "'7.20"
$ws.Range("E26").Value = "'  +2.52%  "
$ws.Range("E27").Value = "'  +1.44%  "
$ws.Range("E28").Value = "'  +1.26%  "
$ws.Range("E29").Value = "'  +0.04%  "
$ws.Range("E30").Value = "'  +0.56%  "
$ws.Range("E31").Value = "'  +2.12%  "
$ws.Range("E32").Value = "'  +2.99%  "
$ws.Range("E33").Value = "'  +3.78%  "
$ws.Range("E34").Value = "'  +1.05%  "
$ws.Range("D35").Value = "'1.446.04"
$ws.Range("B36").Value = "'RenderToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.45"
$ws.Range("E36").Value = "'  +10.15%  "
$ws.Range("B37").Value = "'ImmutableX"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.654"
$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("D38").Value = "'0.0192"
$ws.Range("E38").Value = "'  +3.98%  "
$ws.Range("E39").Value = "'  +0.82%  "
$ws.Range("D40").Value = "'80.64"
$ws.Range("D41").Value = "'2.38"
$ws.Range("E41").Value = "'  +0.82%  "
$ws.Range("E42").Value = "'  +1.73%  "
$ws.Range("E43").Value = "'  +1.20%  "
$ws.Range("D44").Value = "'13.54"
$ws.Range("E44").Value = "'  -0.54%  "
$ws.Range("D45").Value = "'6.08"
$ws.Range("E45").Value = "'  +4.50%  "
$ws.Range("D46").Value = "'0.0508"
$ws.Range("E46").Value = "'  +1.89%  "
$ws.Range("D47").Value = "'1.08"
$ws.Range("E47").Value = "'  -0.22%  "
$ws.Range("E48").Value = "'  -1.79%  "
$ws.Range("D49").Value = "'1.947.88"
$ws.Range("E49").Value = "'  +1.62%  "
$ws.Range("D50").Value = "'106.09"
$ws.Range("E50").Value = "'  -1.34%  "
$ws.Range("E51").Value = "'  +0.03%  "
